$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Spon2"
$ws.Cells.Item(2,3).Value = "Itga5"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.54332
$ws.Cells.Item(2,8).Value = 1.62996
$ws.Cells.Item(2,9).Value = 0.09277509850694737
$ws.Cells.Item(2,10).Value = 0.09480543614915297
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 28.85518433333334
$ws.Cells.Item(2,14).Value = 86.56555300000001
$ws.Cells.Item(2,15).Value = 0.1999651185353207
$ws.Cells.Item(2,16).Value = 0.2044513327926365
$ws.Cells.Item(2,17).Value = 15.67759875198667
$ws.Cells.Item(2,18).Value = 141.09838876788
$ws.Cells.Item(2,19).Value = 0.01855178357006778
$ws.Cells.Item(2,20).Value = 0.01938309777668153

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Spon2"
$ws.Cells.Item(3,3).Value = "Itga5"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.54332
$ws.Cells.Item(3,8).Value = 1.62996
$ws.Cells.Item(3,9).Value = 0.09277509850694737
$ws.Cells.Item(3,10).Value = 0.09480543614915297
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 51.17424933333334
$ws.Cells.Item(3,14).Value = 153.522748
$ws.Cells.Item(3,15).Value = 0.3546352265743414
$ws.Cells.Item(3,16).Value = 0.3625914622481308
$ws.Cells.Item(3,17).Value = 27.80399314778667
$ws.Cells.Item(3,18).Value = 250.23593833008
$ws.Cells.Item(3,19).Value = 0.03290131807946813
$ws.Cells.Item(3,20).Value = 0.03437564172239318

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Spon2"
$ws.Cells.Item(4,3).Value = "Itga5"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.54332
$ws.Cells.Item(4,8).Value = 1.62996
$ws.Cells.Item(4,9).Value = 0.09277509850694737
$ws.Cells.Item(4,10).Value = 0.09480543614915297
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 29.393479
$ws.Cells.Item(4,14).Value = 88.180437
$ws.Cells.Item(4,15).Value = 0.2036954761578358
$ws.Cells.Item(4,16).Value = 0.2082653809291453
$ws.Cells.Item(4,17).Value = 15.97006501028
$ws.Cells.Item(4,18).Value = 143.73058509252
$ws.Cells.Item(4,19).Value = 0.01889786786596277
$ws.Cells.Item(4,20).Value = 0.0197446902737571

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Spon2"
$ws.Cells.Item(5,3).Value = "Itga5"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.54332
$ws.Cells.Item(5,8).Value = 1.62996
$ws.Cells.Item(5,9).Value = 0.09277509850694737
$ws.Cells.Item(5,10).Value = 0.09480543614915297
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 9.499066500000001
$ws.Cells.Item(5,14).Value = 18.998133
$ws.Cells.Item(5,15).Value = 0.0658280999596015
$ws.Cells.Item(5,16).Value = 0.04486996822421697
$ws.Cells.Item(5,17).Value = 5.161032810780001
$ws.Cells.Item(5,18).Value = 30.96619686468001
$ws.Cells.Item(5,19).Value = 0.006107208458277207
$ws.Cells.Item(5,20).Value = 0.004253916907495525

# Row 6: ECs -> Resolving-Mac
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Spon2"
$ws.Cells.Item(6,3).Value = "Itga5"
$ws.Cells.Item(6,4).Value = "Resolving-Mac"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.54332
$ws.Cells.Item(6,8).Value = 1.62996
$ws.Cells.Item(6,9).Value = 0.09277509850694737
$ws.Cells.Item(6,10).Value = 0.09480543614915297
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 25.37910966666666
$ws.Cells.Item(6,14).Value = 76.137329
$ws.Cells.Item(6,15).Value = 0.1758760787729007
$ws.Cells.Item(6,16).Value = 0.1798218558058706
$ws.Cells.Item(6,17).Value = 13.78897786409333
$ws.Cells.Item(6,18).Value = 124.10080077684
$ws.Cells.Item(6,19).Value = 0.0163169205331715
$ws.Cells.Item(6,20).Value = 0.01704808946882566

# Row 7: FAPs -> ECs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Spon2"
$ws.Cells.Item(7,3).Value = "Itga5"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 4.616901666666666
$ws.Cells.Item(7,8).Value = 13.850705
$ws.Cells.Item(7,9).Value = 0.7883632241071366
$ws.Cells.Item(7,10).Value = 0.805616167573593
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 28.85518433333334
$ws.Cells.Item(7,14).Value = 86.56555300000001
$ws.Cells.Item(7,15).Value = 0.1999651185353207
$ws.Cells.Item(7,16).Value = 0.2044513327926365
$ws.Cells.Item(7,17).Value = 133.2215486405406
$ws.Cells.Item(7,18).Value = 1198.993937764865
$ws.Cells.Item(7,19).Value = 0.1576451455574711
$ws.Cells.Item(7,20).Value = 0.1647092991797171

# Row 8: FAPs -> FAPs
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Spon2"
$ws.Cells.Item(8,3).Value = "Itga5"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 4.616901666666666
$ws.Cells.Item(8,8).Value = 13.850705
$ws.Cells.Item(8,9).Value = 0.7883632241071366
$ws.Cells.Item(8,10).Value = 0.805616167573593
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 51.17424933333334
$ws.Cells.Item(8,14).Value = 153.522748
$ws.Cells.Item(8,15).Value = 0.3546352265743414
$ws.Cells.Item(8,16).Value = 0.3625914622481308
$ws.Cells.Item(8,17).Value = 236.2664770374822
$ws.Cells.Item(8,18).Value = 2126.39829333734
$ws.Cells.Item(8,19).Value = 0.2795813706041127
$ws.Cells.Item(8,20).Value = 0.2921095442112443

# Row 9: FAPs -> Inflammatory-Mac
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Spon2"
$ws.Cells.Item(9,3).Value = "Itga5"
$ws.Cells.Item(9,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 4.616901666666666
$ws.Cells.Item(9,8).Value = 13.850705
$ws.Cells.Item(9,9).Value = 0.7883632241071366
$ws.Cells.Item(9,10).Value = 0.805616167573593
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 29.393479
$ws.Cells.Item(9,14).Value = 88.180437
$ws.Cells.Item(9,15).Value = 0.2036954761578358
$ws.Cells.Item(9,16).Value = 0.2082653809291453
$ws.Cells.Item(9,17).Value = 135.7068021842317
$ws.Cells.Item(9,18).Value = 1221.361219658085
$ws.Cells.Item(9,19).Value = 0.1605860223198298
$ws.Cells.Item(9,20).Value = 0.1677819580223925

# Row 10: FAPs -> MuSCs
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Spon2"
$ws.Cells.Item(10,3).Value = "Itga5"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 4.616901666666666
$ws.Cells.Item(10,8).Value = 13.850705
$ws.Cells.Item(10,9).Value = 0.7883632241071366
$ws.Cells.Item(10,10).Value = 0.805616167573593
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 9.499066500000001
$ws.Cells.Item(10,14).Value = 18.998133
$ws.Cells.Item(10,15).Value = 0.0658280999596015
$ws.Cells.Item(10,16).Value = 0.04486996822421697
$ws.Cells.Item(10,17).Value = 43.8562559556275
$ws.Cells.Item(10,18).Value = 263.137535733765
$ws.Cells.Item(10,19).Value = 0.0518964531209983
$ws.Cells.Item(10,20).Value = 0.03614797183994257

# Row 11: FAPs -> Resolving-Mac
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Spon2"
$ws.Cells.Item(11,3).Value = "Itga5"
$ws.Cells.Item(11,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 4.616901666666666
$ws.Cells.Item(11,8).Value = 13.850705
$ws.Cells.Item(11,9).Value = 0.7883632241071366
$ws.Cells.Item(11,10).Value = 0.805616167573593
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 25.37910966666666
$ws.Cells.Item(11,14).Value = 76.137329
$ws.Cells.Item(11,15).Value = 0.1758760787729007
$ws.Cells.Item(11,16).Value = 0.1798218558058706
$ws.Cells.Item(11,17).Value = 117.1728537185494
$ws.Cells.Item(11,18).Value = 1054.555683466945
$ws.Cells.Item(11,19).Value = 0.1386542325047247
$ws.Cells.Item(11,20).Value = 0.1448673943202967

# Row 12: MuSCs -> ECs
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Spon2"
$ws.Cells.Item(12,3).Value = "Itga5"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.073119
$ws.Cells.Item(12,8).Value = 0.219357
$ws.Cells.Item(12,9).Value = 0.01248550104492653
$ws.Cells.Item(12,10).Value = 0.01275874012697842
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 28.85518433333334
$ws.Cells.Item(12,14).Value = 86.56555300000001
$ws.Cells.Item(12,15).Value = 0.1999651185353207
$ws.Cells.Item(12,16).Value = 0.2044513327926365
$ws.Cells.Item(12,17).Value = 2.109862223269
$ws.Cells.Item(12,18).Value = 18.988760009421
$ws.Cells.Item(12,19).Value = 0.002496664696421604
$ws.Cells.Item(12,20).Value = 0.00260854142371563

# Row 13: MuSCs -> FAPs
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Spon2"
$ws.Cells.Item(13,3).Value = "Itga5"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.073119
$ws.Cells.Item(13,8).Value = 0.219357
$ws.Cells.Item(13,9).Value = 0.01248550104492653
$ws.Cells.Item(13,10).Value = 0.01275874012697842
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 51.17424933333334
$ws.Cells.Item(13,14).Value = 153.522748
$ws.Cells.Item(13,15).Value = 0.3546352265743414
$ws.Cells.Item(13,16).Value = 0.3625914622481308
$ws.Cells.Item(13,17).Value = 3.741809937004001
$ws.Cells.Item(13,18).Value = 33.676289433036
$ws.Cells.Item(13,19).Value = 0.004427798491961698
$ws.Cells.Item(13,20).Value = 0.004626210239085007

# Row 14: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(14,1).Value = "MuSCs"
$ws.Cells.Item(14,2).Value = "Spon2"
$ws.Cells.Item(14,3).Value = "Itga5"
$ws.Cells.Item(14,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 0.3333333333333333
$ws.Cells.Item(14,7).Value = 0.073119
$ws.Cells.Item(14,8).Value = 0.219357
$ws.Cells.Item(14,9).Value = 0.01248550104492653
$ws.Cells.Item(14,10).Value = 0.01275874012697842
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 29.393479
$ws.Cells.Item(14,14).Value = 88.180437
$ws.Cells.Item(14,15).Value = 0.2036954761578358
$ws.Cells.Item(14,16).Value = 0.2082653809291453
$ws.Cells.Item(14,17).Value = 2.149221791001
$ws.Cells.Item(14,18).Value = 19.342996119009
$ws.Cells.Item(14,19).Value = 0.002543240080415467
$ws.Cells.Item(14,20).Value = 0.002657203872721132

# Row 15: MuSCs -> MuSCs
$ws.Cells.Item(15,1).Value = "MuSCs"
$ws.Cells.Item(15,2).Value = "Spon2"
$ws.Cells.Item(15,3).Value = "Itga5"
$ws.Cells.Item(15,4).Value = "MuSCs"
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = 0.3333333333333333
$ws.Cells.Item(15,7).Value = 0.073119
$ws.Cells.Item(15,8).Value = 0.219357
$ws.Cells.Item(15,9).Value = 0.01248550104492653
$ws.Cells.Item(15,10).Value = 0.01275874012697842
$ws.Cells.Item(15,11).Value = 2
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 9.499066500000001
$ws.Cells.Item(15,14).Value = 18.998133
$ws.Cells.Item(15,15).Value = 0.0658280999596015
$ws.Cells.Item(15,16).Value = 0.04486996822421697
$ws.Cells.Item(15,17).Value = 0.6945622434135001
$ws.Cells.Item(15,18).Value = 4.167373460481
$ws.Cells.Item(15,19).Value = 0.0008218968108311328
$ws.Cells.Item(15,20).Value = 0.0005724842640785637

# Row 16: MuSCs -> Resolving-Mac
$ws.Cells.Item(16,1).Value = "MuSCs"
$ws.Cells.Item(16,2).Value = "Spon2"
$ws.Cells.Item(16,3).Value = "Itga5"
$ws.Cells.Item(16,4).Value = "Resolving-Mac"
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = 0.3333333333333333
$ws.Cells.Item(16,7).Value = 0.073119
$ws.Cells.Item(16,8).Value = 0.219357
$ws.Cells.Item(16,9).Value = 0.01248550104492653
$ws.Cells.Item(16,10).Value = 0.01275874012697842
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 25.37910966666666
$ws.Cells.Item(16,14).Value = 76.137329
$ws.Cells.Item(16,15).Value = 0.1758760787729007
$ws.Cells.Item(16,16).Value = 0.1798218558058706
$ws.Cells.Item(16,17).Value = 1.855695119717
$ws.Cells.Item(16,18).Value = 16.701256077453
$ws.Cells.Item(16,19).Value = 0.002195900965296632
$ws.Cells.Item(16,20).Value = 0.002294300327378089

# Row 17: Resolving-Mac -> ECs
$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Spon2"
$ws.Cells.Item(17,3).Value = "Itga5"
$ws.Cells.Item(17,4).Value = "ECs"
$ws.Cells.Item(17,5).Value = 2
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 0.3762535
$ws.Cells.Item(17,8).Value = 0.752507
$ws.Cells.Item(17,9).Value = 0.0642475070420447
$ws.Cells.Item(17,10).Value = 0.0437690215344491
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 28.85518433333334
$ws.Cells.Item(17,14).Value = 86.56555300000001
$ws.Cells.Item(17,15).Value = 0.1999651185353207
$ws.Cells.Item(17,16).Value = 0.2044513327926365
$ws.Cells.Item(17,17).Value = 10.85686409856184
$ws.Cells.Item(17,18).Value = 65.14118459137102
$ws.Cells.Item(17,19).Value = 0.01284726036126132
$ws.Cells.Item(17,20).Value = 0.008948634787747727

# Row 18: Resolving-Mac -> FAPs
$ws.Cells.Item(18,1).Value = "Resolving-Mac"
$ws.Cells.Item(18,2).Value = "Spon2"
$ws.Cells.Item(18,3).Value = "Itga5"
$ws.Cells.Item(18,4).Value = "FAPs"
$ws.Cells.Item(18,5).Value = 2
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 0.3762535
$ws.Cells.Item(18,8).Value = 0.752507
$ws.Cells.Item(18,9).Value = 0.0642475070420447
$ws.Cells.Item(18,10).Value = 0.0437690215344491
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 51.17424933333334
$ws.Cells.Item(18,14).Value = 153.522748
$ws.Cells.Item(18,15).Value = 0.3546352265743414
$ws.Cells.Item(18,16).Value = 0.3625914622481308
$ws.Cells.Item(18,17).Value = 19.25449042153933
$ws.Cells.Item(18,18).Value = 115.526942529236
$ws.Cells.Item(18,19).Value = 0.02278442921669212
$ws.Cells.Item(18,20).Value = 0.01587027351934583

# Row 19: Resolving-Mac -> Inflammatory-Mac
$ws.Cells.Item(19,1).Value = "Resolving-Mac"
$ws.Cells.Item(19,2).Value = "Spon2"
$ws.Cells.Item(19,3).Value = "Itga5"
$ws.Cells.Item(19,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19,5).Value = 2
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 0.3762535
$ws.Cells.Item(19,8).Value = 0.752507
$ws.Cells.Item(19,9).Value = 0.0642475070420447
$ws.Cells.Item(19,10).Value = 0.0437690215344491
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 29.393479
$ws.Cells.Item(19,14).Value = 88.180437
$ws.Cells.Item(19,15).Value = 0.2036954761578358
$ws.Cells.Item(19,16).Value = 0.2082653809291453
$ws.Cells.Item(19,17).Value = 11.0593993509265
$ws.Cells.Item(19,18).Value = 66.356396105559
$ws.Cells.Item(19,19).Value = 0.01308692653888321
$ws.Cells.Item(19,20).Value = 0.009115571942768004

# Row 20: Resolving-Mac -> MuSCs
$ws.Cells.Item(20,1).Value = "Resolving-Mac"
$ws.Cells.Item(20,2).Value = "Spon2"
$ws.Cells.Item(20,3).Value = "Itga5"
$ws.Cells.Item(20,4).Value = "MuSCs"
$ws.Cells.Item(20,5).Value = 2
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = 0.3762535
$ws.Cells.Item(20,8).Value = 0.752507
$ws.Cells.Item(20,9).Value = 0.0642475070420447
$ws.Cells.Item(20,10).Value = 0.0437690215344491
$ws.Cells.Item(20,11).Value = 2
$ws.Cells.Item(20,12).Value = 1
$ws.Cells.Item(20,13).Value = 9.499066500000001
$ws.Cells.Item(20,14).Value = 18.998133
$ws.Cells.Item(20,15).Value = 0.0658280999596015
$ws.Cells.Item(20,16).Value = 0.04486996822421697
$ws.Cells.Item(20,17).Value = 3.574057017357751
$ws.Cells.Item(20,18).Value = 14.296228069431
$ws.Cells.Item(20,19).Value = 0.004229291315718919
$ws.Cells.Item(20,20).Value = 0.001963914605455799

# Row 21: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(21,1).Value = "Resolving-Mac"
$ws.Cells.Item(21,2).Value = "Spon2"
$ws.Cells.Item(21,3).Value = "Itga5"
$ws.Cells.Item(21,4).Value = "Resolving-Mac"
$ws.Cells.Item(21,5).Value = 2
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = 0.3762535
$ws.Cells.Item(21,8).Value = 0.752507
$ws.Cells.Item(21,9).Value = 0.0642475070420447
$ws.Cells.Item(21,10).Value = 0.0437690215344491
$ws.Cells.Item(21,11).Value = 3
$ws.Cells.Item(21,12).Value = 1
$ws.Cells.Item(21,13).Value = 25.37910966666666
$ws.Cells.Item(21,14).Value = 76.137329
$ws.Cells.Item(21,15).Value = 0.1758760787729007
$ws.Cells.Item(21,16).Value = 0.1798218558058706
$ws.Cells.Item(21,17).Value = 9.548978838967166
$ws.Cells.Item(21,18).Value = 57.293873033803
$ws.Cells.Item(21,19).Value = 0.01129959960948914
$ws.Cells.Item(21,20).Value = 0.007870626679131753

# Row 22: Inflammatory-Mac -> ECs
$ws.Cells.Item(22,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(22,2).Value = "Spon2"
$ws.Cells.Item(22,3).Value = "Itga5"
$ws.Cells.Item(22,4).Value = "ECs"
$ws.Cells.Item(22,5).Value = 1
$ws.Cells.Item(22,6).Value = 0.3333333333333333
$ws.Cells.Item(22,7).Value = 0.2467186666666667
$ws.Cells.Item(22,8).Value = 0.740156
$ws.Cells.Item(22,9).Value = 0.04212866929894484
$ws.Cells.Item(22,10).Value = 0.04305063461582644
$ws.Cells.Item(22,11).Value = 3
$ws.Cells.Item(22,12).Value = 1
$ws.Cells.Item(22,13).Value = 28.85518433333334
$ws.Cells.Item(22,14).Value = 86.56555300000001
$ws.Cells.Item(22,15).Value = 0.1999651185353207
$ws.Cells.Item(22,16).Value = 0.2044513327926365
$ws.Cells.Item(22,17).Value = 7.11911260514089
$ws.Cells.Item(22,18).Value = 64.07201344626802
$ws.Cells.Item(22,19).Value = 0.008424264350098831
$ws.Cells.Item(22,20).Value = 0.008801759624774528

# Row 23: Inflammatory-Mac -> FAPs
$ws.Cells.Item(23,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(23,2).Value = "Spon2"
$ws.Cells.Item(23,3).Value = "Itga5"
$ws.Cells.Item(23,4).Value = "FAPs"
$ws.Cells.Item(23,5).Value = 1
$ws.Cells.Item(23,6).Value = 0.3333333333333333
$ws.Cells.Item(23,7).Value = 0.2467186666666667
$ws.Cells.Item(23,8).Value = 0.740156
$ws.Cells.Item(23,9).Value = 0.04212866929894484
$ws.Cells.Item(23,10).Value = 0.04305063461582644
$ws.Cells.Item(23,11).Value = 3
$ws.Cells.Item(23,12).Value = 1
$ws.Cells.Item(23,13).Value = 51.17424933333334
$ws.Cells.Item(23,14).Value = 153.522748
$ws.Cells.Item(23,15).Value = 0.3546352265743414
$ws.Cells.Item(23,16).Value = 0.3625914622481308
$ws.Cells.Item(23,17).Value = 12.62564256318756
$ws.Cells.Item(23,18).Value = 113.630783068688
$ws.Cells.Item(23,19).Value = 0.01494031018210681
$ws.Cells.Item(23,20).Value = 0.01560979255606251

# Row 24: Inflammatory-Mac -> Inflammatory-Mac
$ws.Cells.Item(24,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(24,2).Value = "Spon2"
$ws.Cells.Item(24,3).Value = "Itga5"
$ws.Cells.Item(24,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(24,5).Value = 1
$ws.Cells.Item(24,6).Value = 0.3333333333333333
$ws.Cells.Item(24,7).Value = 0.2467186666666667
$ws.Cells.Item(24,8).Value = 0.740156
$ws.Cells.Item(24,9).Value = 0.04212866929894484
$ws.Cells.Item(24,10).Value = 0.04305063461582644
$ws.Cells.Item(24,11).Value = 3
$ws.Cells.Item(24,12).Value = 1
$ws.Cells.Item(24,13).Value = 29.393479
$ws.Cells.Item(24,14).Value = 88.180437
$ws.Cells.Item(24,15).Value = 0.2036954761578358
$ws.Cells.Item(24,16).Value = 0.2082653809291453
$ws.Cells.Item(24,17).Value = 7.251919947574667
$ws.Cells.Item(24,18).Value = 65.267279528172
$ws.Cells.Item(24,19).Value = 0.00858141935274457
$ws.Cells.Item(24,20).Value = 0.00896595681750654

# Row 25: Inflammatory-Mac -> MuSCs
$ws.Cells.Item(25,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(25,2).Value = "Spon2"
$ws.Cells.Item(25,3).Value = "Itga5"
$ws.Cells.Item(25,4).Value = "MuSCs"
$ws.Cells.Item(25,5).Value = 1
$ws.Cells.Item(25,6).Value = 0.3333333333333333
$ws.Cells.Item(25,7).Value = 0.2467186666666667
$ws.Cells.Item(25,8).Value = 0.740156
$ws.Cells.Item(25,9).Value = 0.04212866929894484
$ws.Cells.Item(25,10).Value = 0.04305063461582644
$ws.Cells.Item(25,11).Value = 2
$ws.Cells.Item(25,12).Value = 1
$ws.Cells.Item(25,13).Value = 9.499066500000001
$ws.Cells.Item(25,14).Value = 18.998133
$ws.Cells.Item(25,15).Value = 0.0658280999596015
$ws.Cells.Item(25,16).Value = 0.04486996822421697
$ws.Cells.Item(25,17).Value = 2.343597021458
$ws.Cells.Item(25,18).Value = 14.061582128748
$ws.Cells.Item(25,19).Value = 0.002773250253775936
$ws.Cells.Item(25,20).Value = 0.001931680607244508

# Row 26: Inflammatory-Mac -> Resolving-Mac
$ws.Cells.Item(26,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(26,2).Value = "Spon2"
$ws.Cells.Item(26,3).Value = "Itga5"
$ws.Cells.Item(26,4).Value = "Resolving-Mac"
$ws.Cells.Item(26,5).Value = 1
$ws.Cells.Item(26,6).Value = 0.3333333333333333
$ws.Cells.Item(26,7).Value = 0.2467186666666667
$ws.Cells.Item(26,8).Value = 0.740156
$ws.Cells.Item(26,9).Value = 0.04212866929894484
$ws.Cells.Item(26,10).Value = 0.04305063461582644
$ws.Cells.Item(26,11).Value = 3
$ws.Cells.Item(26,12).Value = 1
$ws.Cells.Item(26,13).Value = 25.37910966666666
$ws.Cells.Item(26,14).Value = 76.137329
$ws.Cells.Item(26,15).Value = 0.1758760787729007
$ws.Cells.Item(26,16).Value = 0.1798218558058706
$ws.Cells.Item(26,17).Value = 6.261500098147111
$ws.Cells.Item(26,18).Value = 56.353500883324
$ws.Cells.Item(26,19).Value = 0.007409425160218704
$ws.Cells.Item(26,20).Value = 0.007741445010238365
